$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 55 was a blank template row; fill it in as book #51
$ws.Range("E55").Value = 51
$ws.Range("F55").Value = "Wish I could Tell You"
$ws.Range("G55").Value = "Durjoy Datta"
$ws.Range("H55").Value = "Fiction"

# Row 56 is a new row for book #52
$ws.Range("E56").Value = 52
$ws.Range("F56").Value = "My Father Baliah"
$ws.Range("G56").Value = "Satyanarayana Y B"
$ws.Range("H56").Value = "Memoir"

# Match formatting of the row above (E54:H54) for both new rows
$ws.Range("E54:H54").Copy()
$ws.Range("E55:H56").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Rows.Item(55).RowHeight = $ws.Rows.Item(54).RowHeight
$ws.Rows.Item(56).RowHeight = $ws.Rows.Item(54).RowHeight

$excel.ActiveWindow.ScrollRow = 43
$ws.Range("F57").Select()
